# Weekly price update: insert a new (most recent) weekly record for
# "Hortaliza, Feria Lagunitas de Puerto Montt - Ciboulette" at row 97,
# pushing the existing rows 97:108 down to 98:109.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 97 (shifts 97:108 -> 98:109, copies
# formatting from the row above/below as Excel normally would, and keeps
# column D's date style on the new row).
$ws.Rows.Item(97).Insert()

# Populate the new row with this week's record.
$ws.Cells.Item(97, 1).Value = 4
$ws.Cells.Item(97, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(97, 3).Value = "Los Lagos"
$ws.Cells.Item(97, 4).Value = 44449
$ws.Cells.Item(97, 5).Value = 10
$ws.Cells.Item(97, 6).Value = 100112039
$ws.Cells.Item(97, 7).Value = "Ciboulette"
$ws.Cells.Item(97, 8).Value = "Sin especificar"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 240
$ws.Cells.Item(97, 11).Value = 4500
$ws.Cells.Item(97, 12).Value = 4500
$ws.Cells.Item(97, 13).Value = 4500
$ws.Cells.Item(97, 14).Value = "$/docena de atados"
$ws.Cells.Item(97, 15).Value = "Región Metropolitana"
$ws.Cells.Item(97, 16).Value = 1500
$ws.Cells.Item(97, 17).Value = 3
$ws.Cells.Item(97, 18).Value = "Hortaliza"
